$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 45, shifting existing rows 45..97 down to 46..98.
$ws.Rows("45").Insert()

# Populate the newly inserted row 45 with the new record.
$ws.Cells.Item(45, 1).Value = 4
$ws.Cells.Item(45, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(45, 3).Value = "Los Lagos"
$ws.Cells.Item(45, 4).Value = (Get-Date -Year 2021 -Month 12 -Day 10 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(45, 5).Value = 10
$ws.Cells.Item(45, 6).Value = 100112022
$ws.Cells.Item(45, 7).Value = "Arveja Verde"
$ws.Cells.Item(45, 8).Value = "Sin especificar"
$ws.Cells.Item(45, 9).Value = "Primera"
$ws.Cells.Item(45, 10).Value = 80
$ws.Cells.Item(45, 11).Value = 20000
$ws.Cells.Item(45, 12).Value = 20000
$ws.Cells.Item(45, 13).Value = 20000
$ws.Cells.Item(45, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(45, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(45, 16).Value = 800
$ws.Cells.Item(45, 17).Value = 25
$ws.Cells.Item(45, 18).Value = "Hortaliza"
